# feat: add 2022-Q4 data
#
# - sheet "总计" (total): the old row "2021-Q2" (row 2) becomes the new
#   "2022-Q4" row (holding value 0.17), and a new row 3 is appended holding
#   the original "2021-Q2" figures (holding value 0.11).
# - a new worksheet "2022-Q4" is inserted between "总计" and "2021-Q2",
#   carrying the quarter's fund-holdings breakdown.

function Set-HeaderStyle($rng) {
    # Mirrors the bold / thin-border / center-top header look already used
    # on the other sheets' header rows.
    $rng.HorizontalAlignment = -4108   # xlCenter
    $rng.VerticalAlignment = -4160     # xlTop
    $rng.Font.Bold = $true
    $rng.Borders.LineStyle = 1         # xlContinuous
}

function Set-TextValue($rng, [string]$val) {
    # Force text storage so numeric-looking strings (fund codes with
    # leading zeros, decimal figures kept as text in the source data)
    # don't get silently re-typed as numbers.
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" sheet: rewrite row 2 as 2022-Q4, append row 3 as 2021-Q2.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.17

$total.Range("A3").Value = 1
Set-HeaderStyle $total.Range("A3")
$total.Range("B3").Value = "2021-Q2"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.11

# ---------------------------------------------------------------------
# 2. Insert a new worksheet "2022-Q4" right before the existing
#    "2021-Q2" sheet, so the tab order becomes 总计, 2022-Q4, 2021-Q2.
# ---------------------------------------------------------------------
$oldQ2 = $wb.Worksheets.Item("2021-Q2")
$q4 = $wb.Worksheets.Add($oldQ2)
$q4.Name = "2022-Q4"

# ---------------------------------------------------------------------
# 3. Populate the new sheet's fund-holdings table.
# ---------------------------------------------------------------------
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"
Set-HeaderStyle $q4.Range("B1:H1")

$q4.Range("A2").Value = 0
Set-HeaderStyle $q4.Range("A2")
Set-TextValue $q4.Range("B2") "014016"
$q4.Range("C2").Value = "中信建投品质优选一年持有期混合A"
Set-TextValue $q4.Range("D2") "4.82"
Set-TextValue $q4.Range("E2") "76.97"
Set-TextValue $q4.Range("F2") "2.79"
Set-TextValue $q4.Range("G2") "0.1345"
$q4.Range("H2").Value = 6

$q4.Range("A3").Value = 1
Set-HeaderStyle $q4.Range("A3")
Set-TextValue $q4.Range("B3") "014017"
$q4.Range("C3").Value = "中信建投品质优选一年持有期混合C"
Set-TextValue $q4.Range("D3") "1.34"
Set-TextValue $q4.Range("E3") "76.97"
Set-TextValue $q4.Range("F3") "2.79"
Set-TextValue $q4.Range("G3") "0.0374"
$q4.Range("H3").Value = 6
